$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Anfernee Simons"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Portland Trail Blazers"

$ws.Range("A3").Value = "Keyonte George"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Utah Jazz"

$ws.Range("A4").Value = "Anthony Edwards"
$ws.Range("B4").Value = "SG,SF"
$ws.Range("C4").Value = "Minnesota Timberwolves"

$ws.Range("A5").Value = "Fred VanVleet"
$ws.Range("B5").Value = "PG"
$ws.Range("C5").Value = "Houston Rockets"

$ws.Range("A6").Value = "Paul George"
$ws.Range("B6").Value = "SG,SF,PF"
$ws.Range("C6").Value = "Philadelphia 76ers"

$ws.Range("A7").Value = "Jaren Jackson Jr."
$ws.Range("B7").Value = "PF,C"
$ws.Range("C7").Value = "Memphis Grizzlies"

$ws.Range("A8").Value = "Jayson Tatum"
$ws.Range("B8").Value = "SF,PF"
$ws.Range("C8").Value = "Boston Celtics"

$ws.Range("A9").Value = "Giannis Antetokounmpo"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Milwaukee Bucks"

$ws.Range("A10").Value = "Yves Missi"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "New Orleans Pelicans"

$ws.Range("A11").Value = "James Harden"
$ws.Range("B11").Value = "PG,SG"
$ws.Range("C11").Value = "LA Clippers"

$ws.Range("A12").Value = "Ivica Zubac"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "LA Clippers"

$ws.Range("A13").Value = "Zion Williamson"
$ws.Range("B13").Value = "PF,C"
$ws.Range("C13").Value = "New Orleans Pelicans"

$ws.Range("A14").Value = "Amen Thompson"
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "Houston Rockets"

$ws.Range("A15").Value = "P.J. Washington"
$ws.Range("B15").Value = "PF"
$ws.Range("C15").Value = "Dallas Mavericks"

$ws.Range("A16").Value = "Jaylin Williams"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Oklahoma City Thunder"

$ws.Range("A17").Value = "Jonathan Kuminga"
$ws.Range("B17").Value = "SF,PF"
$ws.Range("C17").Value = "Golden State Warriors"

$ws.Range("A18").Value = "Goga Bitadze"
$ws.Range("B18").Value = "C"
$ws.Range("C18").Value = "Orlando Magic"

$ws.Range("A19").Value = "Bradley Beal"
$ws.Range("B19").Value = "PG,SG,SF"
$ws.Range("C19").Value = "Phoenix Suns"
